# Updating subjective IQ run demo settings
#
# Clear the demo/sample rows (2-5) back to a blank run template, leaving
# only the date-formatted (but empty) Timestamp cells in column B, and
# move the active selection to A2 (first data row) ready for the next run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2:J5").ClearContents()
$ws.Range("A2").Select() | Out-Null
